$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("medium")

$ws.Cells.Item(8, 2).Value = " @KeyPressed"
$ws.Cells.Item(9, 2).Value = " @KeyPressed"
$ws.Cells.Item(10, 2).Value = " @KeyPressed"
$ws.Cells.Item(11, 2).Value = " @KeyPressed"
$ws.Cells.Item(12, 2).Value = " @KeyPressed"
$ws.Cells.Item(13, 2).Value = " @KeyPressed"
$ws.Cells.Item(14, 2).Value = " @KeyPressed"
$ws.Cells.Item(15, 2).Value = " @KeyPressed"
$ws.Cells.Item(16, 2).Value = " @KeyPressed"
$ws.Cells.Item(17, 2).Value = " @KeyPressed"
$ws.Cells.Item(18, 2).Value = " @KeyPressed"
$ws.Cells.Item(19, 2).Value = " @KeyPressed"
$ws.Cells.Item(20, 1).Value = "What is the probability that the neuron IS behind the jade square? "
$ws.Cells.Item(20, 2).Value = " @seven25ths"
$ws.Cells.Item(21, 2).Value = " @eighteen25ths"
$ws.Cells.Item(21, 4).Value = "."
$ws.Cells.Item(22, 2).Value = " @four25ths"
$ws.Cells.Item(23, 2).Value = " @two9ths"
$ws.Cells.Item(23, 4).Value = "."
$ws.Cells.Item(24, 2).Value = " @KeyPressed"
$ws.Cells.Item(24, 3).Value = "."
$ws.Cells.Item(24, 4).Value = "."
$ws.Cells.Item(25, 2).Value = " @KeyPressed"
$ws.Cells.Item(25, 3).Value = "."
$ws.Cells.Item(25, 4).Value = "."
$ws.Cells.Item(26, 2).Value = " @six25ths"
$ws.Cells.Item(27, 2).Value = " @eighteen25ths"
$ws.Cells.Item(28, 1).Value = "What is the probability that the neuron is not behind the jade square given it is behind the lilac square?"
$ws.Cells.Item(28, 2).Value = " @four6ths"
$ws.Cells.Item(29, 1).Value = "What is the probability that the neuron is behind the lilac square given it is not behind the jade square?"
$ws.Cells.Item(29, 2).Value = " @two9ths"
$ws.Cells.Item(30, 3).Value = "."
$ws.Cells.Item(30, 4).Value = "."

$ws.Activate()
$ws.Range("A35").Select()
